$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(45940,45943,45944,45945,45946,45947,45950,45951,45952,45953,45954,45957,45958,45959,45960,45961,45964,45965,45966,45967,45968,45971,45972,45973,45974,45975,45978,45979,45980,45981)
$scores = @(0.1285659813084559,0.1287382036470091,0.1291497737191831,0.1291170805715238,0.1306651375941179,0.1304588459185503,0.1314351016396507,0.1328840659403057,0.1349019262052822,0.1334518656093679,0.1335783168632589,0.1334961743382388,0.1335629687932902,0.1323912033437438,0.1325974801181503,0.1303342813790535,0.1304443786085342,0.1306046554983352,0.1309971520841812,0.1318909451504794,0.131755389287099,0.1317243901480579,0.1340786938261886,0.1335983523348238,0.1336548873403932,0.133654410503235,0.1336548426369097,0.1329024155868935,0.1329026689066338,0.1329030563368248)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = "SELL"
    $ws.Cells.Item($row, 3).Value = $scores[$i]
}
